$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF" (match style of existing header H1) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-70 for columns I (I0) and J (IF) ---
$iVals = @(6,4,2,4,6,9,8,8,8,9,7,5,8,4,8,9,6,7,6,7,6,7,6,7,5,6,6,6,6,5,9,7,7,9,8,6,5,8,9,6,9,9,5,5,5,7,5,8,9,6,7,9,6,7,6,6,6,9,10,8,6,7,6,6,5,6,3,3,3)
$jVals = @(6,5,4,5,6,9,8,8,9,9,7,6,8,5,8,9,7,7,7,8,6,7,6,7,6,6,6,7,6,6,9,7,7,9,8,6,5,8,9,6,9,9,5,6,6,7,6,8,9,6,8,9,7,8,7,6,7,9,10,8,6,7,6,7,5,6,3,3,3)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
